$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$rows = @{
    1  = @("58÷7=", "45÷4=", "39÷3=", "57÷4=", "78÷7=")
    5  = @("39÷9=", "10÷6=", "59÷8=", "80÷8=", "30÷3=")
    9  = @("72÷4=", "99÷5=", "91÷9=", "97÷3=", "32÷9=")
    13 = @("55÷5=", "78÷8=", "46÷4=", "55÷4=", "12÷8=")
    17 = @("64÷4=", "71÷6=", "14÷6=", "37÷8=", "95÷4=")
}

foreach ($rowIndex in $rows.Keys) {
    $values = $rows[$rowIndex]
    for ($col = 1; $col -le $values.Length; $col++) {
        $t.Cell($rowIndex, $col).Range.Text = $values[$col - 1]
    }
}
